$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '65.647.60'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -3.35%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.477.04'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("E4").Value = '  +0.07%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '581.98'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.15%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '172.33'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -5.49%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.595'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.99%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '3.475.59'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.130'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -7.67%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.83'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.57%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.410'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.92%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.081.19'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.07%  '

$ws.Range("E14").Value = '  +0.24%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '29.84'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -7.16%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.775.67'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.11%  '

$ws.Range("E17").Value = '  -4.03%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.484.60'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.10%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '5.92'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -4.73%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.89'
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '366.47'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -7.43%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.75'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.77%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.537'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '71.97'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0000121'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.00%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.72'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -6.35%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.178'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.84%  '

$ws.Range("E29").Value = '  +0.19%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '24.08'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.70%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '5.76'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -6.34%  '

$ws.Range("E32").Value = '  -3.54%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -8.19%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '7.04'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.31%  '

$ws.Range("E36").Value = '  -1.74%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '29.46'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +11.99%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '158.68'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.885'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.40%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.788.84'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.29%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -12.64%  '

$ws.Range("E43").Value = '  -6.33%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.29'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -6.71%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0686'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.71%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '39.85'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.33%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '24.19'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -8.10%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0287'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.05%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '305.42'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -7.31%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.821'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.40%  '

$ws.Range("E51").Value = '  -4.17%  '
